# Apply the "PO Forecast" sheet addition + header renames described in the diff.
$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" headers on the existing sheets ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet as the 3rd (last) sheet ---
# Duplicate an existing sheet (rather than Worksheets.Add()) so the new
# sheet inherits the same sheetPr/pageMargins/sheetFormatPr as its siblings,
# then wipe its contents and refill them below.
$ws2.Copy($null, $ws2)
$ws3tmp = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3tmp.Name = "PO Forecast"

# Re-fetch the sheet reference by name: after renaming, the old object handle
# stops forwarding writes to the workbook's live sheet.
$ws3 = $wb.Worksheets.Item("PO Forecast")
$ws3.Cells.Clear()

# Copy header formatting (bold, border, centered) from an existing header row.
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Copy the date-number-format style used in column A down the forecast rows.
$ws1.Range("A2").Copy()
$ws3.Range("A2:A13").PasteSpecial(-4122)

# --- 3. Fill in the forecast data rows ---
$ws3.Range("A2").Value = 45501.99999999999
$ws3.Range("B2").Value = 20
$ws3.Range("C2").Value = 19.99999997272108
$ws3.Range("D2").Value = 20.00000002350962

$ws3.Range("A3").Value = 45515.99999999999
$ws3.Range("B3").Value = 20
$ws3.Range("C3").Value = 19.99999997507853
$ws3.Range("D3").Value = 20.00000002455071

$ws3.Range("A4").Value = 45529.99999999999
$ws3.Range("B4").Value = 20
$ws3.Range("C4").Value = 19.99999997452676
$ws3.Range("D4").Value = 20.00000002562393

$ws3.Range("A5").Value = 45536.99999999999
$ws3.Range("B5").Value = 20
$ws3.Range("C5").Value = 19.99999997499289
$ws3.Range("D5").Value = 20.0000000269322

$ws3.Range("A6").Value = 45543.99999999999
$ws3.Range("B6").Value = 20
$ws3.Range("C6").Value = 19.99999996849521
$ws3.Range("D6").Value = 20.00000002927577

$ws3.Range("A7").Value = 45550.99999999999
$ws3.Range("B7").Value = 20
$ws3.Range("C7").Value = 19.99999994374704
$ws3.Range("D7").Value = 20.00000005083889

$ws3.Range("A8").Value = 45557.99999999999
$ws3.Range("B8").Value = 20
$ws3.Range("C8").Value = 19.99999989174547
$ws3.Range("D8").Value = 20.00000010117005

$ws3.Range("A9").Value = 45564.99999999999
$ws3.Range("B9").Value = 20
$ws3.Range("C9").Value = 19.99999982391161
$ws3.Range("D9").Value = 20.0000001750508

$ws3.Range("A10").Value = 45571.99999999999
$ws3.Range("B10").Value = 20
$ws3.Range("C10").Value = 19.99999975044599
$ws3.Range("D10").Value = 20.0000002450666

$ws3.Range("A11").Value = 45578.99999999999
$ws3.Range("B11").Value = 20
$ws3.Range("C11").Value = 19.99999966582343
$ws3.Range("D11").Value = 20.00000032388823

$ws3.Range("A12").Value = 45585.99999999999
$ws3.Range("B12").Value = 20
$ws3.Range("C12").Value = 19.99999958588048
$ws3.Range("D12").Value = 20.00000042430479

$ws3.Range("A13").Value = 45592.99999999999
$ws3.Range("B13").Value = 20
$ws3.Range("C13").Value = 19.99999946896559
$ws3.Range("D13").Value = 20.00000051810837
